$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data block for rows 2..13 (columns A..T)
$data = @(
  @("ECs","Ltb","Tnfrsf1a","ECs",3,1,29.679736,89.039208,0.9444448799234891,0.9444448799234891,3,1,62.49926866666667,187.497806,0.364999568828264,0.364999568828264,1854.961794219739,16694.65614797765,0.3447219739541351,0.3447219739541351),
  @("ECs","Ltb","Tnfrsf1a","sCs",3,1,29.679736,89.039208,0.9444448799234891,0.9444448799234891,3,1,43.717953,131.153859,0.2553155314530077,0.2553155314530077,1297.537303500408,11677.83573150367,0.2411314464457376,0.2411314464457376),
  @("ECs","Ltb","Tnfrsf1a","M2",3,1,29.679736,89.039208,0.9444448799234891,0.9444448799234891,3,1,45.87732766666667,137.631983,0.2679264121734788,0.2679264121734788,1361.626973532163,12254.64276178947,0.2530417281735124,0.2530417281735124),
  @("ECs","Ltb","Tnfrsf1a","FAPs",3,1,29.679736,89.039208,0.9444448799234891,0.9444448799234891,3,1,19.136526,57.409578,0.1117584875452494,0.1117584875452494,567.9670396371361,5111.703356734224,0.1055497313501038,0.1055497313501038),
  @("M2","Ltb","Tnfrsf1a","ECs",3,1,1.732308333333333,5.196925,0.05512413371417655,0.05512413371417656,3,1,62.49926866666667,187.497806,0.364999568828264,0.364999568828264,108.2680039385056,974.4120354465501,0.02012028503770601,0.02012028503770601),
  @("M2","Ltb","Tnfrsf1a","sCs",3,1,1.732308333333333,5.196925,0.05512413371417655,0.05512413371417656,3,1,43.717953,131.153859,0.2553155314530077,0.2553155314530077,75.73297429817501,681.5967686835751,0.01407404749512164,0.01407404749512165),
  @("M2","Ltb","Tnfrsf1a","M2",3,1,1.732308333333333,5.196925,0.05512413371417655,0.05512413371417656,3,1,45.87732766666667,137.631983,0.2679264121734788,0.2679264121734788,79.47367702803056,715.2630932522751,0.01476921137021042,0.01476921137021043),
  @("M2","Ltb","Tnfrsf1a","FAPs",3,1,1.732308333333333,5.196925,0.05512413371417655,0.05512413371417656,3,1,19.136526,57.409578,0.1117584875452494,0.1117584875452494,33.15036346085,298.35327114765,0.006160589811138463,0.006160589811138464),
  @("sCs","Ltb","Tnfrsf1a","ECs",1,0.3333333333333333,0.013544,0.040632,0.0004309863623343461,0.0004309863623343461,3,1,62.49926866666667,187.497806,0.364999568828264,0.364999568828264,0.8464900948213334,7.618410853392001,0.0001573098364228983,0.0001573098364228983),
  @("sCs","Ltb","Tnfrsf1a","sCs",1,0.3333333333333333,0.013544,0.040632,0.0004309863623343461,0.0004309863623343461,3,1,43.717953,131.153859,0.2553155314530077,0.2553155314530077,0.592115955432,5.329043598888001,0.0001100375121483921,0.0001100375121483921),
  @("sCs","Ltb","Tnfrsf1a","M2",1,0.3333333333333333,0.013544,0.040632,0.0004309863623343461,0.0004309863623343461,3,1,45.87732766666667,137.631983,0.2679264121734788,0.2679264121734788,0.6213625259173334,5.592262733256001,0.0001154726297559403,0.0001154726297559403),
  @("sCs","Ltb","Tnfrsf1a","FAPs",1,0.3333333333333333,0.013544,0.040632,0.0004309863623343461,0.0004309863623343461,3,1,19.136526,57.409578,0.1117584875452494,0.1117584875452494,0.259185108144,2.332665973296,0.00004816638400711537,0.00004816638400711537)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $rowNum = $i + 2
  $rowVals = $data[$i]
  for ($c = 0; $c -lt $rowVals.Length; $c++) {
    $colNum = $c + 1
    $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$c]
  }
}
